# Calibration of the needle: sort the data rows (2-12) by column A (time) ascending.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(53613.535026, -0.000023758669731, -0.000019758389067, -0.0000079040205007),
    @(53624.335026, -0.00017156735296, -0.00014172239333, -0.000053008852357),
    @(53635.935027, -0.0003666249, -0.0003011899, -0.0001014434),
    @(53645.803027, -0.0005438759, -0.0004508612, -0.0001518111),
    @(53656.199028, -0.0007193427, -0.0006048907, -0.0002048983),
    @(53666.935028, -0.0008884488, -0.0007603567, -0.0002556789),
    @(53690.999029, -0.0007062664, -0.0005961034, -0.0002109604),
    @(53702.19903, -0.0005284251, -0.000440739, -0.0001587659),
    @(53712.33503, -0.0003495365, -0.0002892978, -0.0001055855),
    @(53723.667031, -0.00015841520936, -0.00013209913367, -0.000053904856423),
    @(53734.531031, -0.000025406350801, -0.00002165626666, -0.000010544177038)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
